# Updated cryptos list with GitHub Actions - refresh Price (D) and Volume(1h) (E) columns
# Price cells are stored as text (values like "61.560.45" aren't valid numbers),
# so NumberFormat is forced to "@" before assigning, then the style is reset back
# to "Normal" so no stray numeric formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.563.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.97%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.391.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.80%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.63"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.47%  "

$ws.Range("E10").Value = "  -0.90%  "

$ws.Range("E11").Value = "  -1.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.972.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.400.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000170"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.603.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.88%  "

$ws.Range("E18").Value = "  +0.16%  "

$ws.Range("E19").Value = "  -0.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.29%  "

$ws.Range("E23").Value = "  -1.01%  "

$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("E25").Value = "  -3.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.191"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.25%  "

$ws.Range("E28").Value = "  +1.21%  "

$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("E30").Value = "  -0.62%  "

$ws.Range("E31").Value = "  +0.32%  "

$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.31"
$ws.Range("D33").Style = "Normal"

$ws.Range("E34").Value = "  -0.82%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "168.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.427.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.30%  "

$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("E39").Value = "  -0.79%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.52"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.781"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.69%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.41%  "

$ws.Range("E44").Value = "  +1.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.462.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.53%  "

$ws.Range("E47").Value = "  -1.32%  "

$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0262"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("E50").Value = "  -5.50%  "

$ws.Range("E51").Value = "  -0.99%  "
